$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B column values (row 2,3,5,6,7,8)
$ws.Range("B2").Value = 0.2860000133514404
$ws.Range("B3").Value = 0.1040999963879585
$ws.Range("B5").Value = 13.39029979705811
$ws.Range("B6").Value = 28.40950012207031
$ws.Range("B7").Value = 5.599500179290771
$ws.Range("B8").Value = 8.887299537658691

# Add new row 9
$ws.Range("A9").Value = "Пастбище"
$ws.Range("B9").Value = 56.67669677734375
$ws.Range("C9").Value = 0.0073
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.2506999969482422
$ws.Range("F9").Value = 0
